# Update Leave Card 12/22/2023 10:59 AM
#
# This script reproduces (via Excel COM automation) the edits recorded in the
# authoritative OOXML diff:
#   - Two new leave-type codes earned/used: rows 62-64 earn 1.25 each (EARNED
#     column C), and a Forced Leave "FL(1-0-0)" is posted on row 65 (1 day
#     taken, dated 12/29/2023 in column K).
#   - A new "2024" year-separator row (row 66) is inserted/labelled, and the
#     monthly period dates for Jan 2024 - Jan 2025 are filled into column A
#     for rows 67-79.
#   - One additional blank row is appended to the bottom of Table1 (which
#     grows from A8:K136 to A8:K137), pushing the old "final" row down and
#     turning what used to be row 136 into an ordinary interior row.
#   - E9/I9 (BALANCE totals) are formulas and recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1) EARNED postings for Sep/Oct/Nov 2023 (rows 62-64), column C (EARNED).
#    Column G ("EARNED ") is a calculated table column that mirrors C and
#    recalculates automatically.
# ---------------------------------------------------------------------
$ws.Range("C62").Value = 1.25
$ws.Range("C63").Value = 1.25
$ws.Range("C64").Value = 1.25

# ---------------------------------------------------------------------
# 2) New year separator "2024" on row 66 (column A), matching the format
#    of the other year-header cells (e.g. A10 = "2020").  The text must
#    stay text (not become the number 2024), so copy the source
#    formatting first and assign the value with a leading apostrophe.
#    (This is populated before B65 below so the new shared-string table
#    entries land in the same order as the reference edit: "2024" then
#    "FL(1-0-0)".)
# ---------------------------------------------------------------------
$ws.Range("A10").Copy()
$ws.Range("A66").PasteSpecial(-4122)
$ws.Range("A66").Value = "'2024"

# ---------------------------------------------------------------------
# 3) Forced Leave taken on row 65 (Dec 2023 period):
#      B65 = "FL(1-0-0)" (particulars)
#      D65 = 1            (Absence/Undertime W/ Pay)
#      K65 = 12/29/2023    (remarks date, serial 45289), styled like the
#                           other "date filed" cells (e.g. K44).
# ---------------------------------------------------------------------
$ws.Range("B65").Value = "FL(1-0-0)"
$ws.Range("D65").Value = 1

$ws.Range("K44").Copy()
$ws.Range("K65").PasteSpecial(-4122)
$ws.Range("K65").Value = 45289

# ---------------------------------------------------------------------
# 4) Monthly period dates for 2024 (Jan 2024 - Jan 2025), rows 67-79,
#    column A. These cells already carry the correct date number format;
#    only the values are being populated.
# ---------------------------------------------------------------------
$ws.Range("A67").Value = 45292
$ws.Range("A68").Value = 45323
$ws.Range("A69").Value = 45352
$ws.Range("A70").Value = 45383
$ws.Range("A71").Value = 45413
$ws.Range("A72").Value = 45444
$ws.Range("A73").Value = 45474
$ws.Range("A74").Value = 45505
$ws.Range("A75").Value = 45536
$ws.Range("A76").Value = 45566
$ws.Range("A77").Value = 45597
$ws.Range("A78").Value = 45627
$ws.Range("A79").Value = 45658

# ---------------------------------------------------------------------
# 5) Append a new trailing row to Table1 (row 137). Excel's AutoExpand
#    would normally do this as soon as the adjacent row is used; here we
#    do it explicitly:
#      - copy the current last row's (136) formatting + formula down to
#        the new row 137,
#      - restore row 136 to the regular interior-row formatting (copied
#        from row 135), since it is no longer the table's final row,
#      - resize Table1 so its reference covers the new row.
# ---------------------------------------------------------------------
$ws.Range("A136:K136").Copy($ws.Range("A137:K137"))
$ws.Range("G137").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

$ws.Range("A135:K135").Copy()
$ws.Range("A136:K136").PasteSpecial(-4122)

$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A8:K137"))

$excel.CutCopyMode = 0
